$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A6").Value = "Klaus"
$ws.Range("C6").Value = "TEST FAILED"
$ws.Range("B6").Value = "Kleber"

$ws.Range("E6").Value = 42932
$ws.Range("E2").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("F6").Value = 99999

$ws.Range("B7").Select()
